$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match results pulled in by the automatic tracker sync (every 3h).
# Each entry: event_id, fecha, jugador_A, jugador_B, pronostico, cuota.
# resultado/profit (G/H) start blank until the match is settled.
$newRows = @(
    @{A="14428735"; B="2025-08-19"; C="Marton Fucsovics"; D="Tallon Griekspoor"; E="Gana Tallon Griekspoor"; F=1.73},
    @{A="14428733"; B="2025-08-19"; C="Hamad Medjedovic"; D="Jacob Fearnley"; E="Gana Jacob Fearnley"; F=2.75},
    @{A="14428731"; B="2025-08-19"; C="Mackenzie McDonald"; D="Luciano Darderi"; E="Gana Luciano Darderi"; F=1.8},
    @{A="14428724"; B="2025-08-19"; C="Mariano Navone"; D="Marcos Giron"; E="Gana Mariano Navone"; F=2.3},
    @{A="14428712"; B="2025-08-19"; C="Sebastián Báez"; D="Pablo Carreño Busta"; E="Gana Sebastián Báez"; F=2.75},
    @{A="14428713"; B="2025-08-19"; C="Mattia Bellucci"; D="Jaume Munar"; E="Gana Mattia Bellucci"; F=2.63},
    @{A="14469463"; B="2025-08-19"; C="Cristina Bucsa"; D="Alycia Parks"; E="Gana Alycia Parks"; F=2.5},
    @{A="14469461"; B="2025-08-19"; C="Leolia Jeanjean"; D="Elisabetta Cocciaretto"; E="Gana Leolia Jeanjean"; F=2.38},
    @{A="14427998"; B="2025-08-19"; C="Eva Lys"; D="Polina Kudermetova"; E="Gana Polina Kudermetova"; F=3.5},
    @{A="14428000"; B="2025-08-19"; C="Solana Sierra"; D="Elsa Jacquemot"; E="Gana Solana Sierra"; F=2.1},
    @{A="14427994"; B="2025-08-19"; C="Ann Li"; D="Iva Jovic"; E="Gana Iva Jovic"; F=1.8},
    @{A="14466758"; B="2025-08-19"; C="August Holmgren"; D="Zachary Svajda"; E="Gana August Holmgren"; F=2.62},
    @{A="14466775"; B="2025-08-19"; C="Jan Choinski"; D="Terence Atmane"; E="Gana Jan Choinski"; F=4.33},
    @{A="14466776"; B="2025-08-19"; C="Santiago Rodriguez Taverna"; D="Mitchell Krueger"; E="Gana Santiago Rodriguez Taverna"; F=3.4},
    @{A="14466774"; B="2025-08-19"; C="Carlos Taberner"; D="Jerome Kym"; E="Gana Carlos Taberner"; F=2.62},
    @{A="14466748"; B="2025-08-19"; C="Luka Pavlovic"; D="Leandro Riedi"; E="Gana Luka Pavlovic"; F=2.62},
    @{A="14466773"; B="2025-08-19"; C="Michael Zheng"; D="Yasutaka Uchiyama"; E="Gana Yasutaka Uchiyama"; F=3.4},
    @{A="14466747"; B="2025-08-19"; C="Alex Bolt"; D="Oliver Crawford"; E="Gana Oliver Crawford"; F=2.5},
    @{A="14466746"; B="2025-08-19"; C="Garrett Johns"; D="Colton Smith"; E="Gana Garrett Johns"; F=5},
    @{A="14466742"; B="2025-08-19"; C="Henrique Rocha"; D="Pierre-Hugues Herbert"; E="Gana Henrique Rocha"; F=2.75},
    @{A="14466768"; B="2025-08-19"; C="Martin Damm Jr"; D="Stefano Travaglia"; E="Gana Stefano Travaglia"; F=2.75},
    @{A="14466752"; B="2025-08-19"; C="Titouan Droguet"; D="Alibek Kachmazov"; E="Gana Alibek Kachmazov"; F=3.25},
    @{A="14466740"; B="2025-08-19"; C="Alex Barrena"; D="Kimmer Coppejans"; E="Gana Alex Barrena"; F=4.5},
    @{A="14466741"; B="2025-08-19"; C="Alexander Blockx"; D="Joao Lucas Reis Da Silva"; E="Gana Joao Lucas Reis Da Silva"; F=5.5},
    @{A="14466738"; B="2025-08-19"; C="Lukas Neumayer"; D="Ignacio Buse"; E="Gana Lukas Neumayer"; F=3},
    @{A="14466783"; B="2025-08-19"; C="Maria Timofeeva"; D="Veronika Erjavec"; E="Gana Veronika Erjavec"; F=2.5},
    @{A="14466836"; B="2025-08-19"; C="Nuria Brancaccio"; D="Priscilla Hon"; E="Gana Nuria Brancaccio"; F=2.5},
    @{A="14466801"; B="2025-08-19"; C="Anca Todoni"; D="Harriet Dart"; E="Gana Harriet Dart"; F=2.75},
    @{A="14466800"; B="2025-08-19"; C="Eva Vedder"; D="Ayana Akli"; E="Gana Ayana Akli"; F=1.91},
    @{A="14466798"; B="2025-08-19"; C="Iryna Shymanovich"; D="Sinja Kraus"; E="Gana Iryna Shymanovich"; F=2.75},
    @{A="14466790"; B="2025-08-19"; C="Maria Carle"; D="Carol Zhao"; E="Gana Carol Zhao"; F=2.62},
    @{A="14466839"; B="2025-08-19"; C="Viktoriya Tomova"; D="Irene Burillo"; E="Gana Irene Burillo"; F=5},
    @{A="14466853"; B="2025-08-19"; C="Ana Sofia Sanchez"; D="Arianne Hartono"; E="Gana Ana Sofia Sanchez"; F=2.2},
    @{A="14466802"; B="2025-08-19"; C="Elizabeth Mandlik"; D="Jessika Ponchet"; E="Gana Jessika Ponchet"; F=1.8},
    @{A="14466855"; B="2025-08-19"; C="Francesca Jones"; D="Viktoria Hruncakova"; E="Gana Viktoria Hruncakova"; F=4},
    @{A="14466854"; B="2025-08-19"; C="Heather Watson"; D="Ekaterine Gorgodze"; E="Gana Ekaterine Gorgodze"; F=2.38},
    @{A="14466844"; B="2025-08-19"; C="Janice Tjen"; D="Varvara Lepchenko"; E="Gana Varvara Lepchenko"; F=3.25},
    @{A="14466820"; B="2025-08-19"; C="Celine Naef"; D="Leyre Romero Gormaz"; E="Gana Leyre Romero Gormaz"; F=2.1}
)

$startRow = 284
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    foreach ($col in @("A", "B", "C", "D", "E")) {
        $cell = $ws.Range("$col$r")
        $cell.NumberFormat = "@"
        $cell.Value = $row[$col]
        $cell.ClearFormats()
    }

    $ws.Range("F$r").Value = $row["F"]
    $ws.Range("G$r").Value = ""
    $ws.Range("H$r").Value = ""
}
